$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect numeric-looking text cells (columns D and G) from Excel's
# automatic number coercion by formatting as Text before assignment,
# then resetting the style afterwards so no extra cell formatting sticks.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "245.13"
$ws.Range("G2").Value = "5"
$ws.Range("G3").Value = "5"
$ws.Range("D4").Value = "5.198"
$ws.Range("G4").Value = "5"
$ws.Range("D5").Value = "0.05747"
$ws.Range("G5").Value = "5"
$ws.Range("D6").Value = "6.467"
$ws.Range("G6").Value = "5"
$ws.Range("D7").Value = "3.231"
$ws.Range("G7").Value = "5"
$ws.Range("D8").Value = "0.8134"
$ws.Range("G8").Value = "5"
$ws.Range("D9").Value = "0.8670"
$ws.Range("G9").Value = "5"
$ws.Range("D10").Value = "0.1378"
$ws.Range("G10").Value = "5"
$ws.Range("D11").Value = "0.06937"
$ws.Range("G11").Value = "5"
$ws.Range("D12").Value = "0.03163"
$ws.Range("G12").Value = "5"
$ws.Range("D13").Value = "0.02966"
$ws.Range("G13").Value = "5"
$ws.Range("D14").Value = "0.09343"
$ws.Range("G14").Value = "5"
$ws.Range("D15").Value = "3.826"
$ws.Range("G15").Value = "5"
$ws.Range("G16").Value = "5"
$ws.Range("D17").Value = "0.04718"
$ws.Range("G17").Value = "5"
$ws.Range("D18").Value = "0.0005997"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("G18").Value = "5"
$ws.Range("D19").Value = "0.006162"
$ws.Range("G19").Value = "5"
$ws.Range("D20").Value = "0.001234"
$ws.Range("G20").Value = "5"
$ws.Range("D21").Value = "0.004104"
$ws.Range("G21").Value = "5"
$ws.Range("D22").Value = "0.00008696"
$ws.Range("G22").Value = "5"
$ws.Range("G23").Value = "5"
$ws.Range("D24").Value = "2.158"
$ws.Range("G24").Value = "5"
$ws.Range("D25").Value = "0.3186"
$ws.Range("G25").Value = "5"
$ws.Range("G26").Value = "5"
$ws.Range("D27").Value = "0.0002327"
$ws.Range("G27").Value = "5"
$ws.Range("G28").Value = "5"
$ws.Range("G29").Value = "5"
$ws.Range("G30").Value = "5"
$ws.Range("G31").Value = "5"
$ws.Range("G32").Value = "5"
$ws.Range("G33").Value = "5"
$ws.Range("G34").Value = "5"
$ws.Range("G35").Value = "5"
$ws.Range("G36").Value = "5"
$ws.Range("G37").Value = "5"
$ws.Range("G38").Value = "5"
$ws.Range("G39").Value = "5"
$ws.Range("D40").Value = "0.03723"
$ws.Range("G40").Value = "5"
$ws.Range("D41").Value = "0.006227"
$ws.Range("G41").Value = "5"
$ws.Range("D42").Value = "0.1052"
$ws.Range("G42").Value = "5"
$ws.Range("D43").Value = "0.002293"
$ws.Range("E43").Value = "42CEJICEJIWorstin24h"
$ws.Range("G43").Value = "5"
$ws.Range("D44").Value = "0.007491"
$ws.Range("G44").Value = "5"
$ws.Range("D45").Value = "0.00005450"
$ws.Range("G45").Value = "5"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("G46").Value = "5"
$ws.Range("D47").Value = "0.4297"
$ws.Range("G47").Value = "5"
$ws.Range("D48").Value = "0.003238"
$ws.Range("G48").Value = "5"
$ws.Range("D49").Value = "0.00002099"
$ws.Range("G49").Value = "5"
$ws.Range("D50").Value = "0.0001999"
$ws.Range("G50").Value = "5"
$ws.Range("G51").Value = "5"

# Reset styles back to Normal so no residual formatting / style index
# differences remain on the cells we touched.
$ws.Range("D2:D51").Style = "Normal"
$ws.Range("G2:G51").Style = "Normal"
